# Apply the "include commentary on some novavax papers" edit to the
# RSV Data 2024 workbook.
#
# Summary of the change (see commit message / diff):
#  - On the "Papers" sheet, insert a new "dose" column (D) between
#    "Drugname" and "Agegroup".
#  - Rework the block of Novavax rows (9-14): reorder them, add a
#    "RSV F nanoparticle vaccine" drug name + dose info, and add/update
#    commentary notes comparing the various Novavax papers.
#  - Add hyperlinks on the four reworked Novavax rows that have a URL.
#  - Make "Papers" the active/selected sheet (it was "Natural Abs"
#    before), and update the remembered selection on both sheets.

$wb = $excel.ActiveWorkbook
$papers = $wb.Worksheets.Item("Papers")
$natAbs = $wb.Worksheets.Item("Natural Abs")

# --- 1. Insert a new "dose" column before the current "Agegroup" column (D) ---
$papers.Columns("D:D").Insert()
$papers.Range("D1").Value = "dose"

# --- 2. Rewrite the Novavax block (rows 9-14) with the new content ---
# Clear the old contents first so nothing stale is left behind.
$papers.Range("A9:J14").ClearContents()

# Row 9: Novavax / Novavax, the original 2017 nanoparticle vaccine paper
$papers.Range("A9").Value = "Novavax"
$papers.Range("B9").Value = "Novavax"
$papers.Range("C9").Value = "RSV F nanoparticle vaccine"
$papers.Range("D9").Value = "multiple"
$papers.Range("E9").Value = "maternal"
$papers.Range("F9").Value = "Abs, Efficacy"
$papers.Range("H9").Value = "Figure 3"
$papers.Range("I9").Value = "this is the first novavax study frm which they all follow"
$papers.Range("J9").Value = "https://www.sciencedirect.com/science/article/pii/S0264410X17306813"

# Row 10: Novavax_preg_ab_JID2019
$papers.Range("A10").Value = "Novavax_preg_ab_JID2019"
$papers.Range("B10").Value = "Novavax"
$papers.Range("C10").Value = "RSV F nanoparticle vaccine"
$papers.Range("D10").Value = "120 μg of RSV F vaccine adsorbed to 0.4 mg of aluminum"
$papers.Range("E10").Value = "maternal, infants"
$papers.Range("F10").Value = "Abs"
$papers.Range("I10").Value = "this shows very little effect of vaccine - and seems to be the same vaccine as Novavax_pref_abeff_NEJM2020"
$papers.Range("J10").Value = "https://academic.oup.com/jid/article/220/11/1802/5546089"

# Row 11: Novavax_pref_abeff_NEJM2020
$papers.Range("A11").Value = "Novavax_pref_abeff_NEJM2020"
$papers.Range("B11").Value = "Novavax"
$papers.Range("C11").Value = "RSV F nanoparticle vaccine"
$papers.Range("D11").Value = "120 μg of RSV F vaccine adsorbed to 0.4 mg of aluminum"
$papers.Range("E11").Value = "maternal, infants"
$papers.Range("F11").Value = "Abs, Efficacy"
$papers.Range("G11").Value = "Prepare"
$papers.Range("H11").Value = "Table S11 (immunogenicity)"
$papers.Range("I11").Value = "this shows strong immunological effect of vaccine and seems to be the same vaccine as Novavax_preg_ab_JID2019"
$papers.Range("J11").Value = "https://www.nejm.org/doi/full/10.1056/NEJMoa1908380"

# Row 12: (analysis derived from the NEJM2020 paper)
$papers.Range("B12").Value = "Novavax"
$papers.Range("C12").Value = "RSV F nanoparticle vaccine"
$papers.Range("D12").Value = "120 μg of RSV F vaccine adsorbed to 0.4 mg of aluminum"
$papers.Range("E12").Value = "maternal, infants"
$papers.Range("F12").Value = "Abs, Efficacy"
$papers.Range("G12").Value = "Prepare"
$papers.Range("I12").Value = "Analysis of data from Novavax_pref_abeff_NEJM2020"
$papers.Range("J12").Value = "https://academic.oup.com/ofid/article/10/1/ofac693/6986211"

# Row 13: Novavax Resolve press release (elderly)
$papers.Range("A13").Value = "Novavax"
$papers.Range("B13").Value = "Novavax"
$papers.Range("E13").Value = "elderly"
$papers.Range("F13").Value = "Efficacy"
$papers.Range("G13").Value = "Resolve"
$papers.Range("J13").Value = "https://ir.novavax.com/press-releases/2016-09-25-Novavax-Announces-Topline-RSV-F-Vaccine-Data-from-Two-Clinical-Trials-in-Older-Adults"

# Row 14: elderly Abs paper
$papers.Range("B14").Value = "Novavax"
$papers.Range("E14").Value = "elderly"
$papers.Range("F14").Value = "Abs"
$papers.Range("J14").Value = "https://immunityageing.biomedcentral.com/articles/10.1186/s12979-017-0090-7"

# --- 3. Hyperlink the URL cells that now carry commentary ---
$papers.Hyperlinks.Add($papers.Range("J11"), "https://www.nejm.org/doi/full/10.1056/NEJMoa1908380")
$papers.Hyperlinks.Add($papers.Range("J9"), "https://www.sciencedirect.com/science/article/pii/S0264410X17306813")
$papers.Hyperlinks.Add($papers.Range("J10"), "https://academic.oup.com/jid/article/220/11/1802/5546089")
$papers.Hyperlinks.Add($papers.Range("J12"), "https://academic.oup.com/ofid/article/10/1/ofac693/6986211")

# --- 4. Best-fit the new "dose" column's neighbour (Agegroup, now column E) ---
$papers.Columns("E:E").AutoFit()

# --- 5. Update sheet selection / active tab: "Papers" becomes the active
#        sheet (it was "Natural Abs" before), remembering the new
#        selections on both sheets. ---
$natAbs.Select()
$natAbs.Range("A17").Select()

$papers.Select()
$papers.Range("C9").Select()
